$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue "D2" '62.854.56'
Set-TextValue "D3" '3.439.08'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue "D5" '577.00'
$ws.Range("E5").Value = '  -1.09%  '
Set-TextValue "D6" '146.34'
Set-TextValue "D7" '3.439.68'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  -0.02%  '
Set-TextValue "D10" '7.74'
$ws.Range("E10").Value = '  +1.31%  '
Set-TextValue "D11" '0.123'
$ws.Range("E11").Value = '  -1.03%  '
$ws.Range("E12").Value = '  +2.51%  '
Set-TextValue "D13" '4.025.86'
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("E14").Value = '  +2.51%  '
Set-TextValue "D15" '28.87'
$ws.Range("E15").Value = '  -1.91%  '
Set-TextValue "D16" '3.426.32'
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("E17").Value = '  -0.92%  '
Set-TextValue "D18" '62.875.54'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("E19").Value = '  +1.90%  '
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("E21").Value = '  -1.23%  '
Set-TextValue "D22" '384.74'
$ws.Range("E22").Value = '  -2.51%  '
Set-TextValue "D23" '0.559'
$ws.Range("E23").Value = '  -0.51%  '
Set-TextValue "D24" '74.33'
$ws.Range("E24").Value = '  -1.48%  '
$ws.Range("E25").Value = '  -0.18%  '
Set-TextValue "D26" '3.580.29'
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("E27").Value = '  -3.85%  '
$ws.Range("E28").Value = '  -6.11%  '
Set-TextValue "D29" '7.53'
$ws.Range("E29").Value = '  -2.34%  '
$ws.Range("E31").Value = '  -1.41%  '
$ws.Range("E32").Value = '  -2.05%  '
$ws.Range("E33").Value = '  -0.08%  '
Set-TextValue "D34" '23.24'
Set-TextValue "D35" '1.30'
$ws.Range("E35").Value = '  -9.41%  '
Set-TextValue "D38" '31.70'
$ws.Range("E38").Value = '  +3.81%  '
Set-TextValue "D39" '1.57'
$ws.Range("E39").Value = '  -0.49%  '
Set-TextValue "D40" '168.63'
$ws.Range("E40").Value = '  +0.35%  '
Set-TextValue "D41" '3.476.23'
$ws.Range("E41").Value = '  +0.00%  '
Set-TextValue "D42" '0.0767'
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("E43").Value = '  -0.58%  '
Set-TextValue "D44" '42.30'
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("E46").Value = '  -1.19%  '
$ws.Range("E47").Value = '  -3.24%  '
Set-TextValue "D48" '2.560.60'
$ws.Range("E48").Value = '  +1.56%  '
$ws.Range("E49").Value = '  +3.82%  '
Set-TextValue "D50" '6.81'
$ws.Range("E50").Value = '  +1.08%  '
Set-TextValue "D51" '22.58'
$ws.Range("E51").Value = '  -4.30%  '
